$d = $word.ActiveDocument

$replacements = @(
    @("114×7=", "437×7="),
    @("563×8=", "906×7="),
    @("727×8=", "599×3="),
    @("664×4=", "768×7="),
    @("756×9=", "262×2="),
    @("911×5=", "977×4="),
    @("886×4=", "557×9="),
    @("963×2=", "604×7="),
    @("493×4=", "173×6="),
    @("184×3=", "978×3="),
    @("394×5=", "404×7="),
    @("822×9=", "531×6="),
    @("445×8=", "290×4="),
    @("976×4=", "882×2="),
    @("159×7=", "974×5="),
    @("251×6=", "207×6="),
    @("465×2=", "203×9="),
    @("816×3=", "201×3="),
    @("313×2=", "970×5="),
    @("294×8=", "255×4="),
    @("146×7=", "522×8="),
    @("622×3=", "759×6="),
    @("794×8=", "999×6="),
    @("959×9=", "946×2="),
    @("615×8=", "663×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
